# Auto-generated edit script
# Adds 2023-09-30 crime data: increments 2023 (column J) totals (and a couple of
# 2018/2022 corrections in the summary sheets) across Citywide Totals, By Neighborhood,
# and per-neighborhood worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 5739
$ws.Range("E3").Value = 7735
$ws.Range("I3").Value = 7494
$ws.Range("J3").Value = 6098
$ws.Range("J4").Value = 1331
$ws.Range("J5").Value = 467
$ws.Range("J6").Value = 7809
$ws.Range("E7").Value = 26014
$ws.Range("I7").Value = 26230
$ws.Range("J7").Value = 21444

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J6").Value = 196
$ws.Range("J7").Value = 300

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 373
$ws.Range("J3").Value = 409
$ws.Range("J6").Value = 455
$ws.Range("J7").Value = 1349

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J2").Value = 131
$ws.Range("J7").Value = 435

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J3").Value = 330
$ws.Range("J7").Value = 992

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("J4").Value = 15
$ws.Range("J5").Value = 10
$ws.Range("J7").Value = 317

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J2").Value = 197
$ws.Range("J6").Value = 194
$ws.Range("J7").Value = 660

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J6").Value = 193
$ws.Range("J7").Value = 543

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("J6").Value = 87
$ws.Range("J7").Value = 337

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J2").Value = 171
$ws.Range("J6").Value = 158
$ws.Range("J7").Value = 628
$ws.Range("J8").Value = 1349
$ws.Range("J15").Value = 236
$ws.Range("J16").Value = 83
$ws.Range("J19").Value = 629
$ws.Range("J20").Value = 448
$ws.Range("J21").Value = 61
$ws.Range("J24").Value = 67
$ws.Range("J29").Value = 1198
$ws.Range("J33").Value = 992
$ws.Range("J36").Value = 294
$ws.Range("J37").Value = 660
$ws.Range("J41").Value = 135
$ws.Range("J42").Value = 894
$ws.Range("J43").Value = 175
$ws.Range("J44").Value = 163
$ws.Range("J48").Value = 254
$ws.Range("J49").Value = 147
$ws.Range("J51").Value = 264
$ws.Range("J52").Value = 538
$ws.Range("J53").Value = 300
$ws.Range("J54").Value = 417
$ws.Range("J55").Value = 289
$ws.Range("J57").Value = 94
$ws.Range("E63").Value = 351
$ws.Range("I63").Value = 243
$ws.Range("J64").Value = 143
$ws.Range("J65").Value = 543
$ws.Range("J67").Value = 811
$ws.Range("J70").Value = 30
$ws.Range("J75").Value = 62
$ws.Range("J76").Value = 317
$ws.Range("J79").Value = 616
$ws.Range("J83").Value = 435
$ws.Range("J84").Value = 182
$ws.Range("J85").Value = 896
$ws.Range("J86").Value = 133
$ws.Range("J89").Value = 282
$ws.Range("J90").Value = 232
$ws.Range("J92").Value = 67
$ws.Range("J95").Value = 317
$ws.Range("J96").Value = 251
$ws.Range("J97").Value = 178
$ws.Range("J99").Value = 337
$ws.Range("E101").Value = 26014
$ws.Range("I101").Value = 26230
$ws.Range("J101").Value = 21444

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J3").Value = 308
$ws.Range("J4").Value = 63
$ws.Range("J6").Value = 215
$ws.Range("J7").Value = 811

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("J2").Value = 58
$ws.Range("J3").Value = 59
$ws.Range("J7").Value = 182

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("J3").Value = 27
$ws.Range("J7").Value = 147

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J3").Value = 81
$ws.Range("J4").Value = 32
$ws.Range("J6").Value = 201
$ws.Range("J7").Value = 417

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 361
$ws.Range("J3").Value = 417
$ws.Range("J7").Value = 1198

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("J2").Value = 41
$ws.Range("J6").Value = 127
$ws.Range("J7").Value = 254

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J2").Value = 155
$ws.Range("J3").Value = 182
$ws.Range("J4").Value = 31
$ws.Range("J6").Value = 238
$ws.Range("J7").Value = 629

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("J3").Value = 38
$ws.Range("J7").Value = 163

$ws = $wb.Worksheets.Item("River North")
$ws.Range("J2").Value = 48
$ws.Range("J6").Value = 176
$ws.Range("J7").Value = 317

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("J2").Value = 45
$ws.Range("J6").Value = 58
$ws.Range("J7").Value = 158

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("J2").Value = 32
$ws.Range("J6").Value = 73
$ws.Range("J7").Value = 135

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J3").Value = 180
$ws.Range("J4").Value = 41
$ws.Range("J6").Value = 463
$ws.Range("J7").Value = 894

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("J6").Value = 144
$ws.Range("J7").Value = 289

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("J6").Value = 15
$ws.Range("J7").Value = 67

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("J2").Value = 75
$ws.Range("J3").Value = 68
$ws.Range("J7").Value = 251

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("J6").Value = 42
$ws.Range("J7").Value = 61

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J3").Value = 212
$ws.Range("J6").Value = 179
$ws.Range("J7").Value = 616

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("J2").Value = 39
$ws.Range("J7").Value = 143

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J2").Value = 123
$ws.Range("J3").Value = 154
$ws.Range("J4").Value = 40
$ws.Range("J6").Value = 121
$ws.Range("J7").Value = 448

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J3").Value = 96
$ws.Range("J7").Value = 294

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J2").Value = 193
$ws.Range("J6").Value = 203
$ws.Range("J7").Value = 628

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J4").Value = 10
$ws.Range("J7").Value = 236

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("J3").Value = 44
$ws.Range("J7").Value = 171

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("J2").Value = 29
$ws.Range("J3").Value = 20
$ws.Range("J7").Value = 178

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("J2").Value = 17
$ws.Range("J7").Value = 67

$ws = $wb.Worksheets.Item("O'Hare")
$ws.Range("J3").Value = 10
$ws.Range("J7").Value = 30

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J2").Value = 89
$ws.Range("J7").Value = 282

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("J4").Value = 71
$ws.Range("J7").Value = 133

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("J6").Value = 14
$ws.Range("J7").Value = 62

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("J4").Value = 8
$ws.Range("J7").Value = 232

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("J2").Value = 63
$ws.Range("J4").Value = 24
$ws.Range("J6").Value = 101
$ws.Range("J7").Value = 264

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("J6").Value = 39
$ws.Range("J7").Value = 94

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("J4").Value = 19
$ws.Range("J7").Value = 175

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J2").Value = 236
$ws.Range("J3").Value = 322
$ws.Range("J7").Value = 896

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J2").Value = 129
$ws.Range("J3").Value = 166
$ws.Range("J6").Value = 215
$ws.Range("J7").Value = 538

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("J3").Value = 6
$ws.Range("J7").Value = 83
